$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Basic Game rubric")

# Update progress values on rows 2 and 3 (Camera / Animations) from 1 -> 2
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 2

# Clear the "klaar op ..." remark cells in column C that belonged to these rows
$ws.Range("C2").Clear()
$ws.Range("C3").Clear()

# Move the active selection to B3, matching the saved selection state
$ws.Activate()
$ws.Range("B3").Select()
